# Update the "Solucion" sheet with the new Pedido -> Salida assignment
# produced by the beam-search / simulated-annealing run.
$wb = $excel.ActiveWorkbook
$wsSolucion = $wb.Worksheets.Item("Solucion")

$nuevasSalidas = @(
    "S051",
    "S008",
    "S062",
    "S029",
    "S050",
    "S063",
    "S042",
    "S009",
    "S053",
    "S033",
    "S034",
    "S018",
    "S001",
    "S025",
    "S017",
    "S039",
    "S058",
    "S044",
    "S007",
    "S038",
    "S032",
    "S005",
    "S037",
    "S016",
    "S061",
    "S040",
    "S036",
    "S031",
    "S056",
    "S026",
    "S030",
    "S004",
    "S024",
    "S019",
    "S035",
    "S014",
    "S002",
    "S021",
    "S013",
    "S047",
    "S045",
    "S006",
    "S012",
    "S048",
    "S022",
    "S059",
    "S064",
    "S049",
    "S043",
    "S055",
    "S011",
    "S028",
    "S010",
    "S027",
    "S023",
    "S046",
    "S054",
    "S060",
    "S003",
    "S020"
)

for ($i = 0; $i -lt $nuevasSalidas.Count; $i++) {
    $fila = $i + 2
    $wsSolucion.Range("B$fila").Value = $nuevasSalidas[$i]
}

# Update the "Metricas" sheet: Z3 max time changes from 388,5 to 388,4
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B4").Value = "388,4"
